# Refresh market-price derived columns (H-N) across all profession sheets.
# Values below mirror a scheduled external price-feed sync; no formulas are
# involved on this workbook -- every H..N cell holds a literal number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 825.75
$ws.Range("I2").Value = 603.3333
$ws.Range("J2").Value = 1493
$ws.Range("K2").Value = 603.3333
$ws.Range("L2").Value = 1493
$ws.Range("M2").Value = -490.3333
$ws.Range("N2").Value = -1719

$ws.Range("H62").Value = 179762.64
$ws.Range("I62").Value = 196739.1
$ws.Range("K62").Value = 196739.1
$ws.Range("M62").Value = -196115.1

$ws.Range("H65").Value = 179762.64
$ws.Range("I65").Value = 196739.1
$ws.Range("K65").Value = 983695.5
$ws.Range("M65").Value = -980575.5

$ws.Range("H70").Value = 7918.684
$ws.Range("J70").Value = 7247.5
$ws.Range("L70").Value = 21742.5
$ws.Range("N70").Value = -22282.5

$ws.Range("H73").Value = 7918.684
$ws.Range("J73").Value = 7247.5
$ws.Range("L73").Value = 21742.5
$ws.Range("N73").Value = -23614.5

$ws.Range("H88").Value = 19999.666
$ws.Range("J88").Value = 27499.5
$ws.Range("L88").Value = 27499.5
$ws.Range("N88").Value = -28311.5

$ws.Range("H91").Value = 19999.666
$ws.Range("J91").Value = 27499.5
$ws.Range("L91").Value = 27499.5
$ws.Range("N91").Value = -30307.5

$ws.Range("H132").Value = 35721244
$ws.Range("I132").Value = 45461330
$ws.Range("K132").Value = 136383990
$ws.Range("M132").Value = -136381460

$ws.Range("H141").Value = 2335.818
$ws.Range("I141").Value = 2201.75
$ws.Range("J141").Value = 2693.3333
$ws.Range("K141").Value = 6605.25
$ws.Range("L141").Value = 8079.999899999999
$ws.Range("M141").Value = -1425.25
$ws.Range("N141").Value = -18439.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3437.6365
$ws.Range("I2").Value = 3437.6365
$ws.Range("K2").Value = 3437.6365
$ws.Range("M2").Value = -3324.6365

$ws.Range("H45").Value = 1596
$ws.Range("I45").Value = 1596
$ws.Range("K45").Value = 1596
$ws.Range("M45").Value = -1219

$ws.Range("H61").Value = 7433.4688
$ws.Range("I61").Value = 6324.9287
$ws.Range("K61").Value = 6324.9287
$ws.Range("M61").Value = -6112.9287

$ws.Range("H76").Value = 77683.14
$ws.Range("J76").Value = 86464.5
$ws.Range("L76").Value = 86464.5
$ws.Range("N76").Value = -87140.5

$ws.Range("H79").Value = 77683.14
$ws.Range("J79").Value = 86464.5
$ws.Range("L79").Value = 86464.5
$ws.Range("N79").Value = -88804.5

$ws.Range("H116").Value = 3437.6365
$ws.Range("I116").Value = 3437.6365
$ws.Range("K116").Value = 3437.6365
$ws.Range("M116").Value = -1143.6365

$ws.Range("H122").Value = 2986.9092
$ws.Range("I122").Value = 3020.25
$ws.Range("J122").Value = 2898
$ws.Range("K122").Value = 9060.75
$ws.Range("L122").Value = 8694
$ws.Range("M122").Value = -6610.75
$ws.Range("N122").Value = -13594

$ws.Range("H136").Value = 7433.4688
$ws.Range("I136").Value = 6324.9287
$ws.Range("K136").Value = 18974.7861
$ws.Range("M136").Value = -16424.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3437.6365
$ws.Range("I3").Value = 3437.6365
$ws.Range("K3").Value = 3437.6365
$ws.Range("M3").Value = -3323.6365

$ws.Range("H20").Value = 1405.25
$ws.Range("I20").Value = 1496.6
$ws.Range("J20").Value = 1253
$ws.Range("K20").Value = 1496.6
$ws.Range("L20").Value = 1253
$ws.Range("M20").Value = -1249.6
$ws.Range("N20").Value = -1747

$ws.Range("H134").Value = 3140.4348
$ws.Range("I134").Value = 3140.4348
$ws.Range("K134").Value = 9421.304400000001
$ws.Range("M134").Value = -6886.304400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 308.875
$ws.Range("I7").Value = 311
$ws.Range("K7").Value = 311
$ws.Range("M7").Value = -198

$ws.Range("H86").Value = 11595.083
$ws.Range("I86").Value = 6682.5
$ws.Range("K86").Value = 6682.5
$ws.Range("M86").Value = -5559.5

$ws.Range("H89").Value = 11595.083
$ws.Range("I89").Value = 6682.5
$ws.Range("K89").Value = 33412.5
$ws.Range("M89").Value = -27796.5

$ws.Range("H105").Value = 3734.6667
$ws.Range("I105").Value = 3569.182
$ws.Range("K105").Value = 3569.182
$ws.Range("M105").Value = -1822.182

$ws.Range("H122").Value = 10700
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 4200
$ws.Range("M122").Value = -1750

$ws.Range("H132").Value = 1960.1666
$ws.Range("I132").Value = 1863.5
$ws.Range("J132").Value = 2153.5
$ws.Range("K132").Value = 5590.5
$ws.Range("L132").Value = 6460.5
$ws.Range("M132").Value = -3060.5
$ws.Range("N132").Value = -11520.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 679.8570999999999
$ws.Range("J12").Value = 474.0909
$ws.Range("L12").Value = 1422.2727
$ws.Range("N12").Value = -1768.2727

$ws.Range("H131").Value = 4389.4326
$ws.Range("I131").Value = 12243.143
$ws.Range("J131").Value = 2556.9
$ws.Range("K131").Value = 36729.429
$ws.Range("L131").Value = 7670.700000000001
$ws.Range("M131").Value = -31689.429
$ws.Range("N131").Value = -17750.7

$ws.Range("H138").Value = 6195.421
$ws.Range("I138").Value = 5392.75
$ws.Range("J138").Value = 7571.4287
$ws.Range("K138").Value = 16178.25
$ws.Range("L138").Value = 22714.2861
$ws.Range("M138").Value = -11038.25
$ws.Range("N138").Value = -32994.2861

$ws.Range("H140").Value = 3533.4656
$ws.Range("I140").Value = 3253.0981
$ws.Range("J140").Value = 5576.143
$ws.Range("K140").Value = 9759.294300000001
$ws.Range("L140").Value = 16728.429
$ws.Range("M140").Value = -4579.294300000001
$ws.Range("N140").Value = -27088.429

$ws.Range("H141").Value = 8549.263000000001
$ws.Range("I141").Value = 7024.222
$ws.Range("K141").Value = 21072.666
$ws.Range("M141").Value = -15892.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 19900
$ws.Range("J92").Value = 19900
$ws.Range("L92").Value = 19900
$ws.Range("N92").Value = -23644

$ws.Range("H102").Value = 51461.05
$ws.Range("I102").Value = 1542.5294
$ws.Range("J102").Value = 334332.66
$ws.Range("K102").Value = 1542.5294
$ws.Range("L102").Value = 334332.66
$ws.Range("M102").Value = 79.4706000000001
$ws.Range("N102").Value = -337576.66

$ws.Range("H113").Value = 2948.75
$ws.Range("I113").Value = 2798
$ws.Range("K113").Value = 2798
$ws.Range("M113").Value = -628

$ws.Range("H122").Value = 111229496
$ws.Range("I122").Value = 143006720
$ws.Range("K122").Value = 429020160
$ws.Range("M122").Value = -429017710

$ws.Range("H126").Value = 93947.82000000001
$ws.Range("J126").Value = 1700
$ws.Range("L126").Value = 5100
$ws.Range("N126").Value = -10040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("K35").Value = 1000
$ws.Range("M35").Value = -664

$ws.Range("H40").Value = 8414.177
$ws.Range("I40").Value = 7178.5835
$ws.Range("K40").Value = 7178.5835
$ws.Range("M40").Value = -7042.5835

$ws.Range("H61").Value = 7220.75
$ws.Range("J61").Value = 6600
$ws.Range("L61").Value = 6600
$ws.Range("N61").Value = -7004

$ws.Range("H82").Value = 6597.4287
$ws.Range("I82").Value = 8990.786
$ws.Range("K82").Value = 8990.786
$ws.Range("M82").Value = -8629.786

$ws.Range("H85").Value = 6597.4287
$ws.Range("I85").Value = 8990.786
$ws.Range("K85").Value = 8990.786
$ws.Range("M85").Value = -7742.786

$ws.Range("H113").Value = 7220.75
$ws.Range("J113").Value = 6600
$ws.Range("L113").Value = 6600
$ws.Range("N113").Value = -10940

$ws.Range("H132").Value = 8694.948
$ws.Range("I132").Value = 8429.261
$ws.Range("J132").Value = 9076.875
$ws.Range("K132").Value = 25287.783
$ws.Range("L132").Value = 27230.625
$ws.Range("M132").Value = -22757.783
$ws.Range("N132").Value = -32290.625

$ws.Range("H136").Value = 5378.6665
$ws.Range("J136").Value = 4799.2
$ws.Range("L136").Value = 14397.6
$ws.Range("N136").Value = -19497.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 9166.666999999999
$ws.Range("J7").Value = 8750
$ws.Range("L7").Value = 8750
$ws.Range("N7").Value = -8976

$ws.Range("H118").Value = 81167.25
$ws.Range("I118").Value = 72388
$ws.Range("K118").Value = 72388
$ws.Range("M118").Value = -70731

$ws.Range("H122").Value = 3513.348
$ws.Range("I122").Value = 3513.348
$ws.Range("K122").Value = 10540.044
$ws.Range("M122").Value = -8090.044

$ws.Range("H126").Value = 3861.2222
$ws.Range("I126").Value = 3760.7273
$ws.Range("J126").Value = 4966.6665
$ws.Range("K126").Value = 11282.1819
$ws.Range("L126").Value = 14899.9995
$ws.Range("M126").Value = -9112.1819
$ws.Range("N126").Value = -19839.9995

